$d = $word.ActiveDocument

# Locate the two paragraphs to remove by their text and delete their
# entire ranges (including the paragraph mark) in one shot.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Fix wrong consumer type by loading") {
        $startPara = $p
    }
    if ($t -eq "Track names changes") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
